$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. '8.00', '1.00'); Excel's Range.Value
# auto-converts such strings to numbers, which would drop formatting like
# trailing zeros. Force text storage via NumberFormat='@', then restore the
# default 'Normal' style so no stray style id is left on the cell.
$dCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D32", "D33", "D34", "D36", "D38", "D42", "D43", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "67.523.73"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "3.508.05"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "604.33"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").Value = "150.08"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("D7").Value = "3.507.37"
$ws.Range("E7").Value = "  -4.48%  "
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("D14").Value = "4.099.61"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "31.56"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "3.494.76"
$ws.Range("E16").Value = "  -4.55%  "
$ws.Range("D17").Value = "67.362.49"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "15.04"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").Value = "447.39"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").Value = "8.99"
$ws.Range("E22").Value = "  -12.81%  "
$ws.Range("D23").Value = "0.621"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").Value = "77.46"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "0.0000130"
$ws.Range("E25").Value = "  +5.91%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "3.646.76"
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -9.05%  "
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("E31").Value = "  -7.01%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "0.165"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "25.74"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").Value = "3.495.28"
$ws.Range("E36").Value = "  -4.90%  "
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("D38").Value = "8.00"
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "174.82"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "0.0876"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -6.94%  "
$ws.Range("D45").Value = "0.881"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").Value = "45.46"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("E47").Value = "  -5.79%  "
$ws.Range("E48").Value = "  +6.10%  "
$ws.Range("D49").Value = "2.56"
$ws.Range("E49").Value = "  -6.08%  "
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -3.94%  "

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }
